$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) value updates, rows 2-93 ---
$ws.Range("C2").Value = 2565.769164761619
$ws.Range("C3").Value = 2771.04675450926
$ws.Range("C4").Value = 6923.341809163824
$ws.Range("C5").Value = 2870.311589353206
$ws.Range("C6").Value = 1909.084588129339
$ws.Range("C7").Value = 6128.19547247793
$ws.Range("C8").Value = 3972.630273980753
$ws.Range("C9").Value = 4729.735976516416
$ws.Range("C10").Value = 752.7964806390332
$ws.Range("C11").Value = 1250.795760575873
$ws.Range("C12").Value = 29601.42505091757
$ws.Range("C13").Value = 3799.441100542036
$ws.Range("C14").Value = 7397.509860835168
$ws.Range("C15").Value = 0
$ws.Range("C18").Value = 492.3430015592067
$ws.Range("C20").Value = 2635.191563620458
$ws.Range("C21").Value = 369.2024078290272
$ws.Range("C22").Value = 1577.487171555845
$ws.Range("C23").Value = 5660.517066940175
$ws.Range("C24").Value = 2024.117324382548
$ws.Range("C25").Value = 4479.398934239905
$ws.Range("C26").Value = 5360.226632400601
$ws.Range("C27").Value = 1410.426304742003
$ws.Range("C28").Value = 1385.890384668919
$ws.Range("C29").Value = 3587.183047009039
$ws.Range("C30").Value = 7453.823475007535
$ws.Range("C31").Value = 0
$ws.Range("C33").Value = 2860.874335573629
$ws.Range("C35").Value = 1223.631935023299
$ws.Range("C36").Value = 389.9389667216314
$ws.Range("C37").Value = 7500.041066630049
$ws.Range("C38").Value = 1657.651524528445
$ws.Range("C39").Value = 5745.422744292303
$ws.Range("C41").Value = 3487.613616731733
$ws.Range("C42").Value = 2634.85005236495
$ws.Range("C43").Value = 7179.116970062444
$ws.Range("C44").Value = 0
$ws.Range("C45").Value = 2887.250212489506
$ws.Range("C47").Value = 1299.811672673934
$ws.Range("C48").Value = 419.1838602515346
$ws.Range("C49").Value = 7563.992777076393
$ws.Range("C50").Value = 1716.389195271215
$ws.Range("C51").Value = 5955.175904294275
$ws.Range("C53").Value = 3405.472039138021
$ws.Range("C54").Value = 2632.058233068435
$ws.Range("C55").Value = 6978.952586250825
$ws.Range("C56").Value = 0
$ws.Range("C57").Value = 2286.013198234259
$ws.Range("C58").Value = 3008.669179463094
$ws.Range("C59").Value = 2648.294169302945
$ws.Range("C60").Value = 449.4203771491282
$ws.Range("C61").Value = 3137.260298393558
$ws.Range("C62").Value = 730.3063521039821
$ws.Range("C63").Value = 3353.623382286602
$ws.Range("C64").Value = 10646.03446486957
$ws.Range("C65").Value = 6753.607115829548
$ws.Range("C66").Value = 558.2093442539386
$ws.Range("C67").Value = 1775.027517189621
$ws.Range("C68").Value = 4861.287098802361
$ws.Range("C69").Value = 5996.49696468919
$ws.Range("C70").Value = 0
$ws.Range("C72").Value = 2361.056581219794
$ws.Range("C73").Value = 3012.536723186288
$ws.Range("C74").Value = 2703.742092148914
$ws.Range("C75").Value = 482.6390663355013
$ws.Range("C76").Value = 3210.869677115934
$ws.Range("C77").Value = 729.1196658666737
$ws.Range("C78").Value = 3305.422815235401
$ws.Range("C79").Value = 10617.47465504905
$ws.Range("C80").Value = 6487.899081675427
$ws.Range("C81").Value = 579.0880693780265
$ws.Range("C82").Value = 1836.014008604312
$ws.Range("C83").Value = 4944.191641077407
$ws.Range("C84").Value = 6114.227214287786
$ws.Range("C85").Value = 0
$ws.Range("C87").Value = 514.0573067519859
$ws.Range("C88").Value = 1875.732161108182
$ws.Range("C89").Value = 6411.986543373589
$ws.Range("C90").Value = 584.2111078769213
$ws.Range("C91").Value = 7476.621011558085
$ws.Range("C93").Value = 0

# --- Column AL (Colony) flips 0 -> 1 ---
$ws.Range("AL23").Value = 1
$ws.Range("AL39").Value = 1
$ws.Range("AL51").Value = 1
$ws.Range("AL62").Value = 1
$ws.Range("AL77").Value = 1
